$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D23").Value = "f"
$ws.Range("E23").Value = "f"

$ws.Range("E23").Select()
